$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" query (cell B2 on the startup sheet) used to return an
# extra `Cohort` column sourced from an OPTIONAL MATCH on (co:cohort).
# The commit drops that trailing column from the Cypher RETURN clause,
# so remove it from the end of the stored query text.
$cell = $ws.Range("B2")
$text = $cell.Value()
$cohortSuffix = ",`n        coalesce(co.cohort_description, '') AS ``Cohort``"

if ($text.EndsWith($cohortSuffix)) {
    $cell.Value = $text.Substring(0, $text.Length - $cohortSuffix.Length)
}

# The author's cursor ended up back on the query cell they just trimmed,
# rather than on D2 (the file name cell) where it was before.
$ws.Range("B2").Select()
